$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the two rows that were dropped from the dataset entirely ---
# Original row 28 = "SC 92", original row 26 = "RM 232".
# Delete the higher-numbered row first so the second delete's row index
# still points at the intended row.
$ws.Rows(28).Delete()
$ws.Rows(26).Delete()

# --- Apply the individual cell value changes (values now at their FINAL
# post-deletion row numbers, 2-33) ---

# Row 2 (RM 2): E2 newly missing
$ws.Range("E2").Value = ""

# Row 5 (RM 14): E5 newly filled in
$ws.Range("E5").Value = -5

# Row 6 (RM 21): D6 and E6 newly filled in
$ws.Range("D6").Value = -14.2
$ws.Range("E6").Value = -5.7

# Row 8 (RM 38): D8 newly missing
$ws.Range("D8").Value = ""

# Row 10 (RM 52 a): E10 newly missing
$ws.Range("E10").Value = ""

# Row 12 (RM 81): D12 newly filled in
$ws.Range("D12").Value = -14.1

# Row 13 (RM 88): E13 newly missing
$ws.Range("E13").Value = ""

# Row 14 (RM 90): D14 newly missing
$ws.Range("D14").Value = ""

# Row 17 (RM 116): D17 newly filled in
$ws.Range("D17").Value = -14.7

# Row 18 (RM 120): D18 newly filled in
$ws.Range("D18").Value = -15.2

# Row 19 (RM 125): D19 newly missing
$ws.Range("D19").Value = ""

# Row 20 (RM 134): D20 newly missing
$ws.Range("D20").Value = ""

# Row 23 (RM 140): D23 newly filled in
$ws.Range("D23").Value = -13.9

# Row 24 (RM 142a): E24 newly filled in
$ws.Range("E24").Value = -8.1

# Row 27 (SC 101, after shift): B27 newly filled in, D27 newly missing
$ws.Range("B27").Value = -20.4
$ws.Range("D27").Value = ""

# Row 28 (SC 105, after shift): E28 newly missing
$ws.Range("E28").Value = ""

# Row 29 (SC 119, after shift): B29 newly missing
$ws.Range("B29").Value = ""

# Row 30 (SC 120, after shift): E30 newly filled in
$ws.Range("E30").Value = -5.7

# Row 32 (SC 193, after shift): B32 newly missing
$ws.Range("B32").Value = ""
